# Generate Report for Handoff
# For the zh-cn and de-de sheets, rows 4-7 (the "Ready for handoff" items)
# get their Priority bumped from "low" to "ht" and their Latest Handoff
# Datetime refreshed, reflecting that a handoff report/xliff was just
# generated for them.

$wb = $excel.ActiveWorkbook

$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

for ($r = 4; $r -le 7; $r++) {
    $zh.Cells.Item($r, 5).Value = "ht"
    $zh.Cells.Item($r, 8).Value = "2016-09-04 02:34:44"

    $de.Cells.Item($r, 5).Value = "ht"
    $de.Cells.Item($r, 8).Value = "2016-09-04 02:34:49"
}
